$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "Fri Oct 06 11:20:35 EDT 2023",
    "Fri Oct 06 11:20:44 EDT 2023",
    "Fri Oct 06 11:20:54 EDT 2023",
    "Fri Oct 06 11:21:04 EDT 2023",
    "Fri Oct 06 11:21:14 EDT 2023",
    "Fri Oct 06 11:21:24 EDT 2023",
    "Fri Oct 06 11:21:34 EDT 2023",
    "Fri Oct 06 11:21:44 EDT 2023",
    "Fri Oct 06 11:21:54 EDT 2023",
    "Fri Oct 06 11:22:04 EDT 2023",
    "Fri Oct 06 11:22:14 EDT 2023",
    "Fri Oct 06 11:22:24 EDT 2023",
    "Fri Oct 06 11:22:35 EDT 2023",
    "Fri Oct 06 11:22:45 EDT 2023",
    "Fri Oct 06 11:22:55 EDT 2023",
    "Fri Oct 06 11:23:05 EDT 2023"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Range("B$row").Value = $newValues[$i]
}
